# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Updates column G ("K") values on the active sheet for rows 2-21.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 3
    3  = 2
    4  = 3
    5  = 4
    6  = 3
    7  = 2
    8  = 4
    9  = 0
    10 = 1
    11 = 1
    12 = 1
    13 = 3
    14 = 4
    15 = 2
    16 = 0
    17 = 6
    18 = 5
    19 = 1
    20 = 1
    21 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
